$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 24079.072
$ws.Cells.Item(8, 9).Value = 125155.625
$ws.Cells.Item(8, 11).Value = 375466.875
$ws.Cells.Item(8, 13).Value = -375327.875
$ws.Cells.Item(31, 8).Value = 2458.6
$ws.Cells.Item(31, 9).Value = 2823.25
$ws.Cells.Item(31, 11).Value = 8469.75
$ws.Cells.Item(31, 13).Value = -8239.75
$ws.Cells.Item(39, 8).Value = 217.64285
$ws.Cells.Item(39, 9).Value = 94.42856999999999
$ws.Cells.Item(39, 10).Value = 340.85715
$ws.Cells.Item(39, 11).Value = 283.28571
$ws.Cells.Item(39, 12).Value = 1022.57145
$ws.Cells.Item(39, 13).Value = 12.71429000000001
$ws.Cells.Item(39, 14).Value = -1614.57145
$ws.Cells.Item(55, 8).Value = 1121.5834
$ws.Cells.Item(55, 9).Value = 424
$ws.Cells.Item(55, 10).Value = 1619.8572
$ws.Cells.Item(55, 11).Value = 424
$ws.Cells.Item(55, 12).Value = 1619.8572
$ws.Cells.Item(55, 13).Value = -210
$ws.Cells.Item(55, 14).Value = -2047.8572
$ws.Cells.Item(86, 8).Value = 3761832.2
$ws.Cells.Item(86, 9).Value = 2979.375
$ws.Cells.Item(86, 10).Value = 8773636
$ws.Cells.Item(86, 11).Value = 2979.375
$ws.Cells.Item(86, 12).Value = 8773636
$ws.Cells.Item(86, 13).Value = -1856.375
$ws.Cells.Item(86, 14).Value = -8775882
$ws.Cells.Item(89, 8).Value = 3761832.2
$ws.Cells.Item(89, 9).Value = 2979.375
$ws.Cells.Item(89, 10).Value = 8773636
$ws.Cells.Item(89, 11).Value = 14896.875
$ws.Cells.Item(89, 12).Value = 43868180
$ws.Cells.Item(89, 13).Value = -9280.875
$ws.Cells.Item(89, 14).Value = -43879412
$ws.Cells.Item(137, 8).Value = 5298.25
$ws.Cells.Item(137, 10).Value = 6957.857
$ws.Cells.Item(137, 12).Value = 20873.571
$ws.Cells.Item(137, 14).Value = -25973.571
$ws.Cells.Item(138, 8).Value = 5994.4053
$ws.Cells.Item(138, 9).Value = 4100.263
$ws.Cells.Item(138, 10).Value = 6594.217
$ws.Cells.Item(138, 11).Value = 12300.789
$ws.Cells.Item(138, 12).Value = 19782.651
$ws.Cells.Item(138, 13).Value = -7160.789000000001
$ws.Cells.Item(138, 14).Value = -30062.651
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2993.5
$ws.Cells.Item(32, 9).Value = 2045.386
$ws.Cells.Item(32, 11).Value = 2045.386
$ws.Cells.Item(32, 13).Value = -1758.386
$ws.Cells.Item(74, 8).Value = 9060.1
$ws.Cells.Item(74, 9).Value = 9421.788
$ws.Cells.Item(74, 10).Value = 7355
$ws.Cells.Item(74, 11).Value = 9421.788
$ws.Cells.Item(74, 12).Value = 7355
$ws.Cells.Item(74, 13).Value = -8547.788
$ws.Cells.Item(74, 14).Value = -9103
$ws.Cells.Item(77, 8).Value = 9060.1
$ws.Cells.Item(77, 9).Value = 9421.788
$ws.Cells.Item(77, 10).Value = 7355
$ws.Cells.Item(77, 11).Value = 47108.94
$ws.Cells.Item(77, 12).Value = 36775
$ws.Cells.Item(77, 13).Value = -42740.94
$ws.Cells.Item(77, 14).Value = -45511
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2378.7917
$ws.Cells.Item(94, 9).Value = 803.4761999999999
$ws.Cells.Item(94, 11).Value = 803.4761999999999
$ws.Cells.Item(94, 13).Value = -352.4761999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 459153.4
$ws.Cells.Item(31, 9).Value = 627380.5600000001
$ws.Cells.Item(31, 10).Value = 10547.667
$ws.Cells.Item(31, 11).Value = 627380.5600000001
$ws.Cells.Item(31, 12).Value = 10547.667
$ws.Cells.Item(31, 13).Value = -627085.5600000001
$ws.Cells.Item(31, 14).Value = -11137.667
$ws.Cells.Item(34, 8).Value = 459153.4
$ws.Cells.Item(34, 9).Value = 627380.5600000001
$ws.Cells.Item(34, 10).Value = 10547.667
$ws.Cells.Item(34, 11).Value = 627380.5600000001
$ws.Cells.Item(34, 12).Value = 10547.667
$ws.Cells.Item(34, 13).Value = -627178.5600000001
$ws.Cells.Item(34, 14).Value = -10951.667
$ws.Cells.Item(35, 8).Value = 718.5454999999999
$ws.Cells.Item(35, 9).Value = 764
$ws.Cells.Item(35, 11).Value = 764
$ws.Cells.Item(35, 13).Value = -470
$ws.Cells.Item(99, 8).Value = 8721.315000000001
$ws.Cells.Item(99, 9).Value = 11417
$ws.Cells.Item(99, 11).Value = 11417
$ws.Cells.Item(99, 13).Value = -9919
$ws.Cells.Item(126, 8).Value = 8721.315000000001
$ws.Cells.Item(126, 9).Value = 11417
$ws.Cells.Item(126, 11).Value = 34251
$ws.Cells.Item(126, 13).Value = -31781
$ws.Cells.Item(134, 8).Value = 226442.06
$ws.Cells.Item(134, 9).Value = 2847.5557
$ws.Cells.Item(134, 11).Value = 8542.667099999999
$ws.Cells.Item(134, 13).Value = -6007.667099999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 54470424
$ws.Cells.Item(4, 9).Value = 41320344
$ws.Cells.Item(4, 10).Value = 120220830
$ws.Cells.Item(4, 11).Value = 123961032
$ws.Cells.Item(4, 12).Value = 360662490
$ws.Cells.Item(4, 13).Value = -123960920
$ws.Cells.Item(4, 14).Value = -360662714
$ws.Cells.Item(7, 8).Value = 1234752.4
$ws.Cells.Item(7, 9).Value = 1543299.2
$ws.Cells.Item(7, 11).Value = 4629897.6
$ws.Cells.Item(7, 13).Value = -4629785.6
$ws.Cells.Item(14, 8).Value = 617.4
$ws.Cells.Item(14, 9).Value = 617.4
$ws.Cells.Item(14, 11).Value = 1852.2
$ws.Cells.Item(14, 13).Value = -1679.2
$ws.Cells.Item(107, 8).Value = 98388.38
$ws.Cells.Item(107, 9).Value = 1176.75
$ws.Cells.Item(107, 10).Value = 158210.92
$ws.Cells.Item(107, 11).Value = 3530.25
$ws.Cells.Item(107, 12).Value = 474632.76
$ws.Cells.Item(107, 13).Value = -1610.25
$ws.Cells.Item(107, 14).Value = -478472.76
$ws.Cells.Item(124, 8).Value = 251109.5
$ws.Cells.Item(124, 9).Value = 334032.66
$ws.Cells.Item(124, 10).Value = 2340
$ws.Cells.Item(124, 11).Value = 1002097.98
$ws.Cells.Item(124, 12).Value = 7020
$ws.Cells.Item(124, 13).Value = -997187.98
$ws.Cells.Item(124, 14).Value = -16840
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 75.066666
$ws.Cells.Item(2, 9).Value = 80.70587999999999
$ws.Cells.Item(2, 10).Value = 67.69231000000001
$ws.Cells.Item(2, 11).Value = 80.70587999999999
$ws.Cells.Item(2, 12).Value = 67.69231000000001
$ws.Cells.Item(2, 13).Value = 32.29412000000001
$ws.Cells.Item(2, 14).Value = -293.69231
$ws.Cells.Item(3, 8).Value = 2005167.2
$ws.Cells.Item(3, 10).Value = 1674749.1
$ws.Cells.Item(3, 12).Value = 1674749.1
$ws.Cells.Item(3, 14).Value = -1674981.1
$ws.Cells.Item(102, 8).Value = 931.7143
$ws.Cells.Item(102, 9).Value = 677.2727
$ws.Cells.Item(102, 11).Value = 677.2727
$ws.Cells.Item(102, 13).Value = 944.7273
$ws.Cells.Item(109, 8).Value = 49250
$ws.Cells.Item(109, 10).Value = 49250
$ws.Cells.Item(109, 12).Value = 49250
$ws.Cells.Item(109, 14).Value = -51330
$ws.Cells.Item(122, 8).Value = 2604.1667
$ws.Cells.Item(122, 9).Value = 1400
$ws.Cells.Item(122, 11).Value = 4200
$ws.Cells.Item(122, 13).Value = -1750
$ws.Cells.Item(128, 8).Value = 72933.336
$ws.Cells.Item(128, 10).Value = 75000
$ws.Cells.Item(128, 12).Value = 75000
$ws.Cells.Item(128, 14).Value = -84960
$ws.Cells.Item(135, 8).Value = 125075000
$ws.Cells.Item(135, 10).Value = 125075000
$ws.Cells.Item(135, 12).Value = 125075000
$ws.Cells.Item(135, 14).Value = -125085140
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1053.0435
$ws.Cells.Item(16, 9).Value = 1055.4546
$ws.Cells.Item(16, 11).Value = 1055.4546
$ws.Cells.Item(16, 13).Value = -885.4546
$ws.Cells.Item(56, 8).Value = 40051
$ws.Cells.Item(56, 9).Value = 40051
$ws.Cells.Item(56, 11).Value = 40051
$ws.Cells.Item(56, 13).Value = -39360
$ws.Cells.Item(100, 8).Value = 1763.3334
$ws.Cells.Item(100, 9).Value = 1722.5
$ws.Cells.Item(100, 10).Value = 1845
$ws.Cells.Item(100, 11).Value = 1722.5
$ws.Cells.Item(100, 12).Value = 1845
$ws.Cells.Item(100, 13).Value = -1181.5
$ws.Cells.Item(100, 14).Value = -2927
$ws.Cells.Item(136, 8).Value = 2904.4614
$ws.Cells.Item(136, 9).Value = 2535.261
$ws.Cells.Item(136, 11).Value = 7605.782999999999
$ws.Cells.Item(136, 13).Value = -5055.782999999999
$ws.Cells.Item(140, 8).Value = 49311
$ws.Cells.Item(140, 10).Value = 49311
$ws.Cells.Item(140, 12).Value = 49311
$ws.Cells.Item(140, 14).Value = -59671
$ws.Cells.Item(141, 8).Value = 50833.168
$ws.Cells.Item(141, 10).Value = 50833.168
$ws.Cells.Item(141, 12).Value = 50833.168
$ws.Cells.Item(141, 14).Value = -61193.168
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 8733.444
$ws.Cells.Item(4, 9).Value = 7489.1665
$ws.Cells.Item(4, 11).Value = 7489.1665
$ws.Cells.Item(4, 13).Value = -7376.1665
$ws.Cells.Item(46, 8).Value = 79885.60000000001
$ws.Cells.Item(46, 10).Value = 79885.60000000001
$ws.Cells.Item(46, 12).Value = 79885.60000000001
$ws.Cells.Item(46, 14).Value = -80347.60000000001
$ws.Cells.Item(108, 8).Value = 90006.39999999999
$ws.Cells.Item(108, 10).Value = 90006.39999999999
$ws.Cells.Item(108, 12).Value = 90006.39999999999
$ws.Cells.Item(108, 14).Value = -97686.39999999999
$ws.Cells.Item(110, 8).Value = 69999.5
$ws.Cells.Item(110, 10).Value = 69999.5
$ws.Cells.Item(110, 12).Value = 69999.5
$ws.Cells.Item(110, 14).Value = -78179.5
$ws.Cells.Item(122, 8).Value = 27780868
$ws.Cells.Item(122, 9).Value = 41668308
$ws.Cells.Item(122, 10).Value = 5987.75
$ws.Cells.Item(122, 11).Value = 125004924
$ws.Cells.Item(122, 12).Value = 17963.25
$ws.Cells.Item(122, 13).Value = -125002474
$ws.Cells.Item(122, 14).Value = -22863.25
$ws.Cells.Item(134, 8).Value = 79885.60000000001
$ws.Cells.Item(134, 10).Value = 79885.60000000001
$ws.Cells.Item(134, 12).Value = 239656.8
$ws.Cells.Item(134, 14).Value = -244726.8
$ws.Cells.Item(136, 8).Value = 136710.28
$ws.Cells.Item(136, 9).Value = 38753.273
$ws.Cells.Item(136, 11).Value = 116259.819
$ws.Cells.Item(136, 13).Value = -113709.819
